# Apply "UCLA / Irvine / Modified Colleges" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: establish new shared-string order: UCLA(24), Plumbing Systems Design(25),
# Sequential Program in HVAC Design(26), Southern(27), Certificate in Global Sustainability(28)

# Row 22 (UCLA / Plumbing Systems Design)
$ws.Range("A22").Value = 22
$ws.Range("B22").Value = "UCLA"
$ws.Range("C1").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "Plumbing Systems Design"

# Row 23 (UCLA / Sequential Program in HVAC Design)
$ws.Range("A23").Value = 23
$ws.Range("B23").Value = "UCLA"
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "Sequential Program in HVAC Design"

# Establish "Southern" as the next new shared string (index 27) before
# "Certificate in Global Sustainability" is introduced, by filling the
# whole E17:E32 block now (values beyond row 23 will be cleared back out below).
$ws.Range("E17:E32").Value = "Southern"

# Row 24 (UCLA / Certificate in Global Sustainability)
$ws.Range("A24").Value = 24
$ws.Range("B24").Value = "UCLA"
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = "Certificate in Global Sustainability"

# Row 25 (UCLA / Mechanical Engineering)
$ws.Range("A25").Value = 25
$ws.Range("B25").Value = "UCLA"
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = "Mechanical Engineering"

# Row 26 (UCLA / Electrical Engineering)
$ws.Range("A26").Value = 26
$ws.Range("B26").Value = "UCLA"
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = "Electrical Engineering"

# Row 27 (UCLA / Computer Science)
$ws.Range("A27").Value = 27
$ws.Range("B27").Value = "UCLA"
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = "Computer Science"

# Row 28 (UCLA / Others)
$ws.Range("A28").Value = 28
$ws.Range("B28").Value = "UCLA"
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = "Others"

# --- Step 2: fill column E ("Region") for rows 10-16 with "Northern" and
# clear the over-filled "Southern" values for rows 24-32 (keeping the
# touched row spans at 1:5 without leaving stray empty cells).
$ws.Range("E10:E16").Value = "Northern"
$ws.Range("E24:E32").Value = ""

# --- Step 3: update the saved selection to match the authored state
$ws.Range("B17").Select()
$excel.CutCopyMode = $false
